$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1885723333333333
$ws.Range("H2").Value = 0.565717
$ws.Range("M2").Value = 0.3478976666666666
$ws.Range("N2").Value = 1.043693
$ws.Range("O2").Value = 0.172630997073507
$ws.Range("P2").Value = 0.172630997073507
$ws.Range("Q2").Value = 0.06560387476455555
$ws.Range("R2").Value = 0.5904348728809999
$ws.Range("S2").Value = 0.172630997073507
$ws.Range("T2").Value = 0.172630997073507

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1885723333333333
$ws.Range("H3").Value = 0.565717
$ws.Range("M3").Value = 1.358150333333333
$ws.Range("N3").Value = 4.074451
$ws.Range("O3").Value = 0.6739304936002712
$ws.Range("P3").Value = 0.6739304936002711
$ws.Range("Q3").Value = 0.2561095773741111
$ws.Range("R3").Value = 2.304986196367
$ws.Range("S3").Value = 0.6739304936002712
$ws.Range("T3").Value = 0.6739304936002711

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1885723333333333
$ws.Range("H4").Value = 0.565717
$ws.Range("M4").Value = 0.3092196666666667
$ws.Range("N4").Value = 0.927659
$ws.Range("O4").Value = 0.1534385093262219
$ws.Range("P4").Value = 0.1534385093262219
$ws.Range("Q4").Value = 0.05831027405588889
$ws.Range("R4").Value = 0.524792466503
$ws.Range("S4").Value = 0.1534385093262219
$ws.Range("T4").Value = 0.1534385093262219
